$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.059.77'
$ws.Range('E2').Value = '  -0.81%  '
$ws.Range('D3').Value = '1.808.86'
$ws.Range('E3').Value = '  -2.14%  '
$ws.Range('E4').Value = '  +0.58%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '232.44'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.610'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('E7').Value = '  +0.57%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '40.27'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -5.72%  '
$ws.Range('E9').Value = '  +4.89%  '
$ws.Range('E10').Value = '  -0.74%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0996'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.79%  '
$ws.Range('D12').Value = '2.072.34'
$ws.Range('E12').Value = '  -2.08%  '
$ws.Range('D13').Value = '1.814.15'
$ws.Range('E13').Value = '  -1.75%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.663'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.17%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '11.02'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.94%  '
$ws.Range('E16').Value = '  -2.05%  '
$ws.Range('D17').Value = '35.029.45'
$ws.Range('E17').Value = '  -0.75%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.69'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.52%  '
$ws.Range('D19').Value = '0.0₃0789'
$ws.Range('E19').Value = '  -0.31%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '237.81'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.88'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.51%  '
$ws.Range('E22').Value = '  -1.33%  '
$ws.Range('E23').Value = '  +0.58%  '
$ws.Range('E24').Value = '  +2.84%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '172.04'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.37%  '
$ws.Range('E26').Value = '  -1.46%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.49'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.14%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.120'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.84%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.58'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +19.15%  '
$ws.Range('E30').Value = '  +0.51%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.13'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.94%  '
$ws.Range('E32').Value = '  +3.88%  '
$ws.Range('E33').Value = '  -0.45%  '
$ws.Range('E34').Value = '  -6.24%  '
$ws.Range('E35').Value = '  +5.43%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '91.86'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.81%  '
$ws.Range('E37').Value = '  +0.29%  '
$ws.Range('E38').Value = '  -0.50%  '
$ws.Range('D39').Value = '1.312.29'
$ws.Range('E39').Value = '  -1.78%  '
$ws.Range('E40').Value = '  +0.87%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.997'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.99%  '
$ws.Range('E42').Value = '  +0.35%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '14.48'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.84%  '
$ws.Range('E44').Value = '  -4.96%  '
$ws.Range('E45').Value = '  -2.33%  '
$ws.Range('E46').Value = '  +5.17%  '
$ws.Range('E47').Value = '  -1.60%  '
$ws.Range('D48').Value = '1.988.95'
$ws.Range('E48').Value = '  -1.40%  '
$ws.Range('E49').Value = '  +0.53%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0656'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.16%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '99.37'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.93%  '
